# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# Cell text that merely *looks* numeric (e.g. "10.05", "1.00") must stay as
# literal text, matching the workbook's original inline-string cells. Writing it
# straight through .Value would make Excel's input parser coerce it to a number,
# so each value is entered with a leading apostrophe (Excel's "force text" marker)
# and then ClearFormats() strips the resulting quote-prefix styling, leaving the
# cell back on the default (unstyled) format exactly like its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$updates = @{
    'D2' = '64.267.06'
    'E2' = '  +1.24%  '
    'D3' = '3.149.93'
    'E3' = '  +1.71%  '
    'E4' = '  +0.04%  '
    'D5' = '591.65'
    'E5' = '  +1.58%  '
    'D6' = '147.99'
    'E6' = '  +2.22%  '
    'D8' = '3.139.58'
    'E8' = '  +1.66%  '
    'E9' = '  +1.24%  '
    'E10' = '  +2.31%  '
    'E11' = '  +5.74%  '
    'E12' = '  +0.54%  '
    'D13' = '0.0000248'
    'E13' = '  +0.96%  '
    'D14' = '37.57'
    'E14' = '  +0.94%  '
    'D15' = '3.666.69'
    'E15' = '  +1.71%  '
    'E16' = '  -0.17%  '
    'E17' = '  +2.62%  '
    'D18' = '64.031.48'
    'E18' = '  +1.12%  '
    'D19' = '3.142.54'
    'E19' = '  +1.77%  '
    'D20' = '469.52'
    'E20' = '  +1.67%  '
    'D21' = '14.46'
    'E21' = '  +1.67%  '
    'D22' = '0.733'
    'E22' = '  +1.42%  '
    'D23' = '7.61'
    'E23' = '  +2.05%  '
    'E24' = '  +13.11%  '
    'D25' = '13.17'
    'E25' = '  +2.13%  '
    'D26' = '81.15'
    'E26' = '  -0.14%  '
    'B27' = 'RenderToken'
    'C27' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D27' = '10.05'
    'E27' = '  +11.75%  '
    'B28' = 'Dai'
    'C28' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D28' = '1.00'
    'E28' = '  +0.12%  '
    'E29' = '  +1.84%  '
    'B30' = 'NEARProtocol'
    'C30' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D30' = '7.31'
    'E30' = '  +6.49%  '
    'B31' = 'ImmutableX'
    'C31' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D31' = '2.22'
    'E31' = '  +1.22%  '
    'B32' = 'FirstDigitalUSD'
    'C32' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D32' = '1.00'
    'E32' = '  -0.07%  '
    'E33' = '  +4.11%  '
    'D34' = '27.76'
    'E34' = '  +4.31%  '
    'D35' = '0.0₃0862'
    'E35' = '  +1.30%  '
    'E36' = '  +3.27%  '
    'D37' = '6.18'
    'E37' = '  +3.24%  '
    'E38' = '  +0.29%  '
    'D39' = '3.31'
    'E39' = '  -2.64%  '
    'D40' = '464.68'
    'E40' = '  +7.11%  '
    'D41' = '51.32'
    'E41' = '  +2.05%  '
    'D42' = '9.33'
    'E42' = '  +7.21%  '
    'D43' = '0.292'
    'E43' = '  +8.39%  '
    'D44' = '0.0375'
    'E44' = '  +2.01%  '
    'D45' = '2.892.59'
    'E45' = '  +0.41%  '
    'D46' = '39.83'
    'E46' = '  +10.83%  '
    'E47' = '  +0.41%  '
    'E48' = '  +8.23%  '
    'E49' = '  -0.01%  '
    'E50' = '  +1.10%  '
    'E51' = '  +4.03%  '
}

foreach ($cellRef in $updates.Keys) {
    Set-TextValue $cellRef $updates[$cellRef]
}
